$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# List of (cell, new value) updates derived from the commit diff.
$updates = @(
    @{Cell='D2'; Value='31.495.63'}
    @{Cell='E2'; Value='  +3.80%  '}
    @{Cell='D3'; Value='1.984.98'}
    @{Cell='E3'; Value='  +5.90%  '}
    @{Cell='D4'; Value='0.9908'}
    @{Cell='E4'; Value='  -1.11%  '}
    @{Cell='D5'; Value='0.8342'}
    @{Cell='E5'; Value='  +77.65%  '}
    @{Cell='D6'; Value='252.09'}
    @{Cell='E6'; Value='  +3.38%  '}
    @{Cell='D7'; Value='0.9924'}
    @{Cell='E7'; Value='  -1.01%  '}
    @{Cell='D8'; Value='0.3446'}
    @{Cell='E8'; Value='  +20.23%  '}
    @{Cell='D9'; Value='25.81'}
    @{Cell='E9'; Value='  +17.80%  '}
    @{Cell='D10'; Value='0.06959'}
    @{Cell='E10'; Value='  +8.43%  '}
    @{Cell='D11'; Value='0.8421'}
    @{Cell='E11'; Value='  +17.14%  '}
    @{Cell='D12'; Value='0.08104'}
    @{Cell='E12'; Value='  +4.03%  '}
    @{Cell='D13'; Value='102.27'}
    @{Cell='E13'; Value='  +7.79%  '}
    @{Cell='D14'; Value='1.979.96'}
    @{Cell='E14'; Value='  +5.51%  '}
    @{Cell='D15'; Value='5.500'}
    @{Cell='E15'; Value='  +7.11%  '}
    @{Cell='D16'; Value='275.12'}
    @{Cell='E16'; Value='  -1.21%  '}
    @{Cell='D17'; Value='31.501.85'}
    @{Cell='D18'; Value='14.01'}
    @{Cell='E18'; Value='  +8.32%  '}
    @{Cell='D19'; Value='0.000007892'}
    @{Cell='E19'; Value='  +6.78%  '}
    @{Cell='D20'; Value='2.246.72'}
    @{Cell='E20'; Value='  +4.86%  '}
    @{Cell='D21'; Value='5.675'}
    @{Cell='E21'; Value='  +8.85%  '}
    @{Cell='D22'; Value='0.9953'}
    @{Cell='E22'; Value='  -0.62%  '}
    @{Cell='D23'; Value='0.9890'}
    @{Cell='E23'; Value='  -1.36%  '}
    @{Cell='D24'; Value='6.887'}
    @{Cell='E24'; Value='  +10.49%  '}
    @{Cell='B25'; Value='Stellar'}
    @{Cell='C25'; Value='https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'}
    @{Cell='D25'; Value='0.1619'}
    @{Cell='E25'; Value='  +68.82%  '}
    @{Cell='B26'; Value='Cosmos'}
    @{Cell='C26'; Value='https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'}
    @{Cell='D26'; Value='9.668'}
    @{Cell='E26'; Value='  +7.37%  '}
    @{Cell='D27'; Value='165.93'}
    @{Cell='E27'; Value='  +1.28%  '}
    @{Cell='D28'; Value='19.71'}
    @{Cell='E28'; Value='  +5.70%  '}
    @{Cell='D29'; Value='2.237'}
    @{Cell='E29'; Value='  +19.32%  '}
    @{Cell='D30'; Value='1.557'}
    @{Cell='E30'; Value='  +6.24%  '}
    @{Cell='D31'; Value='1.346'}
    @{Cell='E31'; Value='  -0.40%  '}
    @{Cell='D32'; Value='4.556'}
    @{Cell='E32'; Value='  +8.41%  '}
    @{Cell='E33'; Value='  +6.18%  '}
    @{Cell='D34'; Value='0.05201'}
    @{Cell='E34'; Value='  +8.05%  '}
    @{Cell='D35'; Value='1.219'}
    @{Cell='E35'; Value='  +9.26%  '}
    @{Cell='D36'; Value='0.7446'}
    @{Cell='E36'; Value='  +8.79%  '}
    @{Cell='D37'; Value='2.757'}
    @{Cell='E37'; Value='  +1.64%  '}
    @{Cell='D38'; Value='0.9920'}
    @{Cell='E38'; Value='  -0.94%  '}
    @{Cell='D39'; Value='0.01987'}
    @{Cell='E39'; Value='  +6.35%  '}
    @{Cell='D40'; Value='2.901'}
    @{Cell='E40'; Value='  +3.32%  '}
    @{Cell='D41'; Value='6.586'}
    @{Cell='E41'; Value='  +5.76%  '}
    @{Cell='D42'; Value='78.43'}
    @{Cell='E42'; Value='  +5.79%  '}
    @{Cell='D43'; Value='0.4661'}
    @{Cell='E43'; Value='  +10.68%  '}
    @{Cell='D44'; Value='2.082'}
    @{Cell='E44'; Value='  +7.66%  '}
    @{Cell='D45'; Value='0.8546'}
    @{Cell='E45'; Value='  +3.85%  '}
    @{Cell='D46'; Value='105.51'}
    @{Cell='E46'; Value='  +4.79%  '}
    @{Cell='D47'; Value='0.9931'}
    @{Cell='E47'; Value='  -0.87%  '}
    @{Cell='D48'; Value='9.939'}
    @{Cell='E48'; Value='  +3.75%  '}
    @{Cell='D49'; Value='7.512'}
    @{Cell='E49'; Value='  +9.04%  '}
    @{Cell='D50'; Value='36.51'}
    @{Cell='E50'; Value='  +4.41%  '}
    @{Cell='D51'; Value='0.4268'}
    @{Cell='E51'; Value='  +9.73%  '}
)

foreach ($u in $updates) {
    $c = $ws.Range($u.Cell)
    # Preserve the existing cell style while forcing the value to be
    # written as text so that numeric-looking strings (e.g. "252.09")
    # are not silently converted into Excel numbers.
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $u.Value
    $c.Style = $origStyle
}

Write-Host "Applied $($updates.Count) cell updates"
